$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 (shifts existing rows 35-93 down to 36-94,
# carrying all their data/formatting with them).
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new record.
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(35, 3).Value = 'Ñuble'
$ws.Cells.Item(35, 4).Value = 44804
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = 100112031
$ws.Cells.Item(35, 7).Value = 'Poroto verde'
$ws.Cells.Item(35, 8).Value = 'Magnum'
$ws.Cells.Item(35, 9).Value = 'Primera'
$ws.Cells.Item(35, 10).Value = 60
$ws.Cells.Item(35, 11).Value = 35000
$ws.Cells.Item(35, 12).Value = 35000
$ws.Cells.Item(35, 13).Value = 35000
$ws.Cells.Item(35, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(35, 15).Value = 'Perú'
$ws.Cells.Item(35, 16).Value = 1400
$ws.Cells.Item(35, 17).Value = 25
$ws.Cells.Item(35, 18).Value = 'Hortaliza'

# Match the date cell format used by the other date cells in column D.
$ws.Cells.Item(35, 4).NumberFormat = $ws.Cells.Item(36, 4).NumberFormat()
